$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalculate "Precio Total" (column F) as Cantidad (D) * Precio Unitario (E)
# for each data row, so the total price updates automatically.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 6 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cantidad = $ws.Cells.Item($row, 4).Value2
    $precioUnitario = $ws.Cells.Item($row, 5).Value2
    if ($cantidad -ne $null -and $precioUnitario -ne $null) {
        $ws.Cells.Item($row, 6).Value = $cantidad * $precioUnitario
    }
}
